$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("egresados")

# K4 previously referenced the "adsadsads" string; it now refers to the
# "Ingeniería civil" career text (new shared string, added first).
$ws.Range("K4").Value = "Ingeniería civil"

# F3 previously held a date value; it now becomes free-form text
# (new shared string, added second).
$ws.Range("F3").Value = "asdadsadsasd"

# Selection / active cell moved from G7 to F10.
$ws.Range("F10").Select()
